# AccountCreateJmapping.xlsx update
# - createAccount sheet: flip the ExecutionFlag (column C) from "No" to "Yes"
#   for every account row that hadn't already been flipped.
# - JMapping sheet: update the QA log message / charge description.
# - Leave the JMapping sheet as the active tab (lease create date work),
#   and leave createAccount scrolled/selected near the bottom where work
#   left off.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("createAccount")
$ws2 = $wb.Worksheets.Item("JMapping")

# Rows on createAccount whose ExecutionFlag (column C) needs to become "Yes"
$rows = @(2,3,4,5,6,7,8,10,11,12,13,15,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,39,41,42,43,45,47,48,49,50,51,52,53,54,55,56)

foreach ($r in $rows) {
    $ws1.Cells.Item($r, 3).Value = "Yes"
}

# JMapping sheet content updates (row 2)
$ws2.Range("A2").Value = "QA Testing"
$ws2.Range("C2").Value = "FASB - Charge"

# Restore createAccount's selection to where editing stopped, then
# switch focus to JMapping so it becomes the active/selected tab.
$ws1.Range("A43").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("C2").Select() | Out-Null
